# Refresh the crypto price snapshot (Price / Volume(1h) / Hora columns)
# for the updated scrape run. Values are prefixed with a literal leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# inline-string cells) instead of auto-converting them to numbers/percentages,
# which would silently drop things like trailing zeros ("43.60") or render
# tiny values in scientific notation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.23"
$ws.Range("E2").Value = "'-0.90%"
$ws.Range("G2").Value = "'16"
$ws.Range("D3").Value = "'43.57"
$ws.Range("E3").Value = "'4.85%"
$ws.Range("G3").Value = "'16"
$ws.Range("D4").Value = "'5.488"
$ws.Range("E4").Value = "'-3.39%"
$ws.Range("G4").Value = "'16"
$ws.Range("D5").Value = "'0.08161"
$ws.Range("E5").Value = "'-2.23%"
$ws.Range("G5").Value = "'16"
$ws.Range("D6").Value = "'8.741"
$ws.Range("E6").Value = "'-0.53%"
$ws.Range("G6").Value = "'16"
$ws.Range("D7").Value = "'4.349"
$ws.Range("E7").Value = "'-3.67%"
$ws.Range("G7").Value = "'16"
$ws.Range("D8").Value = "'1.885"
$ws.Range("E8").Value = "'-6.21%"
$ws.Range("G8").Value = "'16"
$ws.Range("D9").Value = "'2.838"
$ws.Range("E9").Value = "'-3.53%"
$ws.Range("G9").Value = "'16"
$ws.Range("D10").Value = "'0.9416"
$ws.Range("E10").Value = "'1.63%"
$ws.Range("G10").Value = "'16"
$ws.Range("D11").Value = "'0.1203"
$ws.Range("E11").Value = "'-7.25%"
$ws.Range("G11").Value = "'16"
$ws.Range("D12").Value = "'0.1898"
$ws.Range("E12").Value = "'-4.26%"
$ws.Range("G12").Value = "'16"
$ws.Range("D13").Value = "'0.09727"
$ws.Range("E13").Value = "'2.36%"
$ws.Range("G13").Value = "'16"
$ws.Range("D14").Value = "'0.04162"
$ws.Range("E14").Value = "'7.12%"
$ws.Range("G14").Value = "'16"
$ws.Range("E15").Value = "'0.79%"
$ws.Range("G15").Value = "'16"
$ws.Range("D16").Value = "'0.001289"
$ws.Range("E16").Value = "'-1.14%"
$ws.Range("G16").Value = "'16"
$ws.Range("D17").Value = "'0.006088"
$ws.Range("E17").Value = "'-0.25%"
$ws.Range("G17").Value = "'16"
$ws.Range("D18").Value = "'3.534"
$ws.Range("E18").Value = "'2.81%"
$ws.Range("G18").Value = "'16"
$ws.Range("E19").Value = "'-0.08%"
$ws.Range("G19").Value = "'16"
$ws.Range("D20").Value = "'8.767"
$ws.Range("E20").Value = "'7.38%"
$ws.Range("G20").Value = "'16"
$ws.Range("D21").Value = "'0.1369"
$ws.Range("E21").Value = "'-0.29%"
$ws.Range("G21").Value = "'16"
$ws.Range("D22").Value = "'0.2499"
$ws.Range("E22").Value = "'-0.47%"
$ws.Range("G22").Value = "'16"
$ws.Range("D23").Value = "'0.04378"
$ws.Range("E23").Value = "'-0.96%"
$ws.Range("G23").Value = "'16"
$ws.Range("D24").Value = "'0.001241"
$ws.Range("E24").Value = "'-2.50%"
$ws.Range("G24").Value = "'16"
$ws.Range("D25").Value = "'0.004314"
$ws.Range("E25").Value = "'-2.13%"
$ws.Range("G25").Value = "'16"
$ws.Range("E26").Value = "'3.00%"
$ws.Range("G26").Value = "'16"
$ws.Range("D27").Value = "'0.0004015"
$ws.Range("E27").Value = "'31.86%"
$ws.Range("G27").Value = "'16"
$ws.Range("G28").Value = "'16"
$ws.Range("G29").Value = "'16"
$ws.Range("G30").Value = "'16"
$ws.Range("G31").Value = "'16"
$ws.Range("G32").Value = "'16"
$ws.Range("G33").Value = "'16"
$ws.Range("G34").Value = "'16"
$ws.Range("G35").Value = "'16"
$ws.Range("G36").Value = "'16"
$ws.Range("G37").Value = "'16"
$ws.Range("G38").Value = "'16"
$ws.Range("D39").Value = "'0.02737"
$ws.Range("E39").Value = "'-2.35%"
$ws.Range("G39").Value = "'16"
$ws.Range("D40").Value = "'0.05722"
$ws.Range("E40").Value = "'2.62%"
$ws.Range("G40").Value = "'16"
$ws.Range("D41").Value = "'0.007904"
$ws.Range("E41").Value = "'1.46%"
$ws.Range("G41").Value = "'16"
$ws.Range("D42").Value = "'0.009754"
$ws.Range("E42").Value = "'4.80%"
$ws.Range("G42").Value = "'16"
$ws.Range("E43").Value = "'-1.53%"
$ws.Range("G43").Value = "'16"
$ws.Range("D44").Value = "'0.002107"
$ws.Range("E44").Value = "'-2.49%"
$ws.Range("G44").Value = "'16"
$ws.Range("D45").Value = "'0.009688"
$ws.Range("E45").Value = "'-7.66%"
$ws.Range("G45").Value = "'16"
$ws.Range("D46").Value = "'0.00007317"
$ws.Range("E46").Value = "'4.63%"
$ws.Range("G46").Value = "'16"
$ws.Range("E47").Value = "'0.55%"
$ws.Range("G47").Value = "'16"
$ws.Range("D48").Value = "'0.003452"
$ws.Range("E48").Value = "'-2.17%"
$ws.Range("G48").Value = "'16"
$ws.Range("E49").Value = "'0.20%"
$ws.Range("G49").Value = "'16"
$ws.Range("E50").Value = "'0.55%"
$ws.Range("G50").Value = "'16"
$ws.Range("E51").Value = "'0.55%"
$ws.Range("G51").Value = "'16"
